$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 242.5
$ws.Range("I58").Value = 242.5
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 727.5
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -577.5
$ws.Range("N58").ClearContents()
$ws.Range("H100").Value = 3308
$ws.Range("I100").Value = 3199.6
$ws.Range("K100").Value = 3199.6
$ws.Range("M100").Value = -2658.6
$ws.Range("H111").Value = 643.6667
$ws.Range("I111").Value = 499.5
$ws.Range("K111").Value = 1498.5
$ws.Range("M111").Value = 1568.5
$ws.Range("H112").Value = 3808
$ws.Range("J112").Value = 4737.25
$ws.Range("L112").Value = 14211.75
$ws.Range("N112").Value = -16427.75
$ws.Range("H116").Value = 37333
$ws.Range("I116").Value = 6999
$ws.Range("K116").Value = 6999
$ws.Range("M116").Value = -3557
$ws.Range("H138").Value = 4745.1816
$ws.Range("J138").Value = 5133.3335
$ws.Range("L138").Value = 15400.0005
$ws.Range("N138").Value = -25680.0005

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 22048.518
$ws.Range("I32").Value = 17755.852
$ws.Range("K32").Value = 17755.852
$ws.Range("M32").Value = -17468.852
$ws.Range("H110").Value = 4188.75
$ws.Range("I110").Value = 1010
$ws.Range("J110").Value = 4642.857
$ws.Range("K110").Value = 1010
$ws.Range("L110").Value = 4642.857
$ws.Range("M110").Value = 1035
$ws.Range("N110").Value = -8732.857

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 13333
$ws.Range("J16").Value = 20000
$ws.Range("L16").Value = 20000
$ws.Range("N16").Value = -20574
$ws.Range("H31").Value = 7104.6665
$ws.Range("I31").Value = 5999
$ws.Range("K31").Value = 5999
$ws.Range("M31").Value = -5704
$ws.Range("H34").Value = 7104.6665
$ws.Range("I34").Value = 5999
$ws.Range("K34").Value = 5999
$ws.Range("M34").Value = -5797
$ws.Range("H58").Value = 5634.6875
$ws.Range("I58").Value = 4046.1428
$ws.Range("K58").Value = 4046.1428
$ws.Range("M58").Value = -3843.1428
$ws.Range("H59").Value = 28819.076
$ws.Range("J59").Value = 34998.332
$ws.Range("L59").Value = 34998.332
$ws.Range("N59").Value = -37288.332
$ws.Range("H60").Value = 7422
$ws.Range("J60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("N60").ClearContents()
$ws.Range("H113").Value = 13333
$ws.Range("J113").Value = 20000
$ws.Range("L113").Value = 20000
$ws.Range("N113").Value = -24340
$ws.Range("H132").Value = 1739.8
$ws.Range("I132").Value = 1674.75
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 5024.25
$ws.Range("L132").Value = 6000
$ws.Range("M132").Value = -2494.25
$ws.Range("N132").Value = -11060
$ws.Range("H134").Value = 14498.5
$ws.Range("I134").Value = 14000
$ws.Range("K134").Value = 42000
$ws.Range("M134").Value = -39465
$ws.Range("H136").Value = 5634.6875
$ws.Range("I136").Value = 4046.1428
$ws.Range("K136").Value = 12138.4284
$ws.Range("M136").Value = -9588.428400000001

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 322.9091
$ws.Range("J2").Value = 415.25
$ws.Range("L2").Value = 2491.5
$ws.Range("N2").Value = -2717.5
$ws.Range("H44").Value = 952.25
$ws.Range("I44").Value = 603
$ws.Range("K44").Value = 1809
$ws.Range("M44").Value = -1411
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
$ws.Range("H86").Value = 850
$ws.Range("H89").Value = 850
$ws.Range("H92").Value = 864.6
$ws.Range("I92").Value = 949.6667
$ws.Range("J92").Value = 737
$ws.Range("K92").Value = 2849.0001
$ws.Range("L92").Value = 2211
$ws.Range("M92").Value = -1601.0001
$ws.Range("N92").Value = -4707

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4999.5
$ws.Range("J80").Value = 4999.5
$ws.Range("L80").Value = 4999.5
$ws.Range("N80").Value = -6995.5
$ws.Range("H83").Value = 4999.5
$ws.Range("J83").Value = 4999.5
$ws.Range("L83").Value = 24997.5
$ws.Range("N83").Value = -34981.5
$ws.Range("H107").Value = 285.25
$ws.Range("I107").Value = 250
$ws.Range("K107").Value = 250
$ws.Range("M107").Value = 1670
$ws.Range("H113").Value = 1169.3334
$ws.Range("I113").Value = 1133.2
$ws.Range("K113").Value = 1133.2
$ws.Range("M113").Value = 1036.8

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H39").Value = 25000
$ws.Range("J39").Value = 25000
$ws.Range("L39").Value = 25000
$ws.Range("N39").Value = -25920
$ws.Range("H64").Value = 49999.5
$ws.Range("J64").Value = 49999.5
$ws.Range("L64").Value = 49999.5
$ws.Range("N64").Value = -50449.5
$ws.Range("H67").Value = 49999.5
$ws.Range("J67").Value = 49999.5
$ws.Range("L67").Value = 49999.5
$ws.Range("N67").Value = -51559.5
$ws.Range("H93").Value = 1697.5
$ws.Range("I93").Value = 1697.5
$ws.Range("K93").Value = 1697.5
$ws.Range("M93").Value = -449.5
$ws.Range("H132").Value = 14050.143
$ws.Range("I132").Value = 14789.643
$ws.Range("J132").Value = 12571.143
$ws.Range("K132").Value = 44368.929
$ws.Range("L132").Value = 37713.429
$ws.Range("M132").Value = -41838.929
$ws.Range("N132").Value = -42773.429
$ws.Range("H136").Value = 6625.875
$ws.Range("I136").Value = 2834.3333
$ws.Range("K136").Value = 8502.999899999999
$ws.Range("M136").Value = -5952.999899999999

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1475.5454
$ws.Range("I113").Value = 801.1177
$ws.Range("K113").Value = 2403.3531
$ws.Range("M113").Value = -233.3531000000003
$ws.Range("H132").Value = 2908.375
$ws.Range("I132").Value = 2908.375
$ws.Range("K132").Value = 8725.125
$ws.Range("M132").Value = -6195.125
$ws.Range("H136").Value = 2545.6667
$ws.Range("I136").Value = 2613.4546
$ws.Range("K136").Value = 7840.3638
$ws.Range("M136").Value = -5290.3638
